# Adds "SchemeType" and "SchemeYear" rows to the config table, just above
# the existing "ResolutionType" row (old row 26, which becomes row 28 once
# the two new rows are inserted above it). Existing rows 26-30 therefore
# shift down to rows 28-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows before the current row 26 ("ResolutionType").
$ws.Rows.Item(26).Resize(2).Insert() | Out-Null

# New row 26: SchemeType
$ws.Cells.Item(26, 1).Value = "SchemeType"
# New row 27: SchemeYear
$ws.Cells.Item(27, 1).Value = "SchemeYear"

# New row 26: B/C
$ws.Cells.Item(26, 2).Value = "Not Scheme Specific"
# New row 27: B/C
$ws.Cells.Item(27, 2).Value = "Not Year Specific"

$ws.Cells.Item(26, 3).Value = "Scheme Type for customer notification"
$ws.Cells.Item(27, 3).Value = "Scheme Year for customer notification"

# Match the formatting used by the other single-row entries in the table
# (column A: left/center/indent; columns B & C: same plus wrap text), and
# the standard 30pt row height used throughout this block.
$ws.Range("A26:C27").HorizontalAlignment = -4131 ; # xlLeft
$ws.Range("A26:C27").VerticalAlignment = -4108 ; # xlCenter
$ws.Range("A26:C27").IndentLevel = 1
$ws.Range("B26:C27").WrapText = $true
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 30

# Grow Table1 so it covers the two newly inserted rows (A1:D30 -> A1:D32);
# this updates both the table ref and its autoFilter range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D32"))

# Reflect the authored selection / scroll position.
$ws.Range("B28").Select()
$excel.ActiveWindow.ScrollRow = 21
